$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = 2.3
$ws.Range("H4").Value = 2.87
$ws.Range("I4").Value = 3.35
$ws.Range("K4").Value = 5.5
$ws.Range("L4").Value = 1.5
$ws.Range("M4").Value = 2.42
$ws.Range("N4").Value = 2.42
$ws.Range("O4").Value = 1.5
$ws.Range("T4").Value = 6.2
$ws.Range("U4").Value = 10.25
$ws.Range("V4").Value = 9.25
$ws.Range("W4").Value = 24
$ws.Range("X4").Value = 22
$ws.Range("Z4").Value = 5.5
$ws.Range("AA4").Value = 5.6
$ws.Range("AB4").Value = 16.5
$ws.Range("AD4").Value = 7.6
$ws.Range("AE4").Value = 16
$ws.Range("AF4").Value = 12
$ws.Range("AG4").Value = 50
$ws.Range("AH4").Value = 37
$ws.Range("G5").Value = 1.75
$ws.Range("Z5").Value = 11
$ws.Range("AG5").Value = 51
$ws.Range("I7").Value = 3.75
$ws.Range("S7").Value = 1.88
$ws.Range("AB7").Value = 15
$ws.Range("AC7").Value = 70
$ws.Range("AD7").Value = 11.25
$ws.Range("AF7").Value = 12.5
$ws.Range("N8").Value = 1.93
$ws.Range("O8").Value = 1.88
$ws.Range("H10").Value = 3.85
$ws.Range("O10").Value = 2.35
$ws.Range("R10").Value = 1.45
$ws.Range("S10").Value = 2.37
$ws.Range("T10").Value = 11.25
$ws.Range("U10").Value = 12
$ws.Range("X10").Value = 13
$ws.Range("Y10").Value = 18
$ws.Range("AH10").Value = 25
$ws.Range("H11").Value = 4.85
$ws.Range("I11").Value = 6.2
$ws.Range("L11").Value = 1.1
$ws.Range("N11").Value = 1.34
$ws.Range("O11").Value = 3
$ws.Range("Q11").Value = 4
$ws.Range("T11").Value = 12
$ws.Range("U11").Value = 9.75
$ws.Range("W11").Value = 11.25
$ws.Range("AB11").Value = 14.5
$ws.Range("AC11").Value = 40
$ws.Range("AD11").Value = 29
$ws.Range("AE11").Value = 50
$ws.Range("AF11").Value = 20
$ws.Range("AI11").Value = 37
$ws.Range("G12").Value = 2.18
$ws.Range("H12").Value = 3
$ws.Range("I12").Value = 3.25
$ws.Range("K12").Value = 6.3
$ws.Range("L12").Value = 1.39
$ws.Range("M12").Value = 2.77
$ws.Range("P12").Value = 1.5
$ws.Range("Q12").Value = 2.4
$ws.Range("R12").Value = 1.85
$ws.Range("S12").Value = 1.85
$ws.Range("T12").Value = 7
$ws.Range("U12").Value = 10.5
$ws.Range("V12").Value = 8.75
$ws.Range("W12").Value = 22
$ws.Range("Z12").Value = 6.3
$ws.Range("AA12").Value = 5.9
$ws.Range("AB12").Value = 14.5
$ws.Range("AC12").Value = 75
$ws.Range("AD12").Value = 8.5
$ws.Range("AE12").Value = 16.5
$ws.Range("AF12").Value = 11.75
$ws.Range("AG12").Value = 45
$ws.Range("AJ12").Value = 700
$ws.Range("G17").Value = 2.8
$ws.Range("I17").Value = 2.55
$ws.Range("L17").Value = 1.4
$ws.Range("M17").Value = 2.75
$ws.Range("T17").Value = 7.5
$ws.Range("U17").Value = 13
$ws.Range("W17").Value = 29
$ws.Range("X17").Value = 26
$ws.Range("AE17").Value = 11
$ws.Range("AF17").Value = 10
